$wb = $excel.ActiveWorkbook

# Helper: write a value as TEXT (shared-string) even when it looks like a
# number, without leaving a permanent custom number-format behind.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- Restricciones_del_follower: regenerated problem coefficients ---
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

Set-TextValue $wsFollower.Range("A2") "5.95 - y"
Set-TextValue $wsFollower.Range("B2") "-5.95"
Set-TextValue $wsFollower.Range("D2") "0.37"
Set-TextValue $wsFollower.Range("E2") "0.8999999999999999"
Set-TextValue $wsFollower.Range("F2") "0"

Set-TextValue $wsFollower.Range("A3") "-0.6000000000000005 - x + y"
Set-TextValue $wsFollower.Range("B3") "-2.3999999999999995"
Set-TextValue $wsFollower.Range("D3") "0.44"
Set-TextValue $wsFollower.Range("E3") "-6.3"
Set-TextValue $wsFollower.Range("F3") "-7.0"

Set-TextValue $wsFollower.Range("A4") "-17.25 + x + 2y"
Set-TextValue $wsFollower.Range("B4") "5.25"
Set-TextValue $wsFollower.Range("D4") "0.0"
Set-TextValue $wsFollower.Range("E4") "-3.4000000000000004"
Set-TextValue $wsFollower.Range("F4") "-3.7"

Set-TextValue $wsFollower.Range("A5") "-16.35 + 4x - y"
Set-TextValue $wsFollower.Range("B5") "3.4499999999999993"
Set-TextValue $wsFollower.Range("D5") "0.07"
Set-TextValue $wsFollower.Range("E5") "2.0"
Set-TextValue $wsFollower.Range("F5") "0"

# --- Punto_modificado: A2 (x) and B2 (y) updated ---
$wsPunto = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $wsPunto.Range("A2") "5.35"
Set-TextValue $wsPunto.Range("B2") "5.95"

# --- Vector_bf: A2 updated ---
# NOTE: Worksheets.Item(name) is case-insensitive, and this workbook has two
# sheets whose names differ only by case ("Vector_bf" vs "Vector_BF"), so we
# must address them positionally (1-based index, in tab order) instead.
$wsVecbf = $wb.Worksheets.Item(5)   # "Vector_bf"
Set-TextValue $wsVecbf.Range("A2") "-1.0"

# --- Vector_BF: A2 and A3 updated ---
$wsVecBF = $wb.Worksheets.Item(6)   # "Vector_BF"
Set-TextValue $wsVecBF.Range("A2") "-9.899999999999999"
Set-TextValue $wsVecBF.Range("A3") "19.0"
